$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 19.5703125
$ws.Columns.Item(2).ColumnWidth = 43.7109375
$ws.Columns.Item(3).ColumnWidth = 28
$ws.Columns.Item(4).ColumnWidth = 28.7109375
$ws.Columns.Item(5).ColumnWidth = 23.140625
$ws.Columns.Item(6).ColumnWidth = 123.85546875

# --- Header row (row 1) ---
$ws.Range("A1").Value = "PART NO."
$ws.Range("B1").Value = "PART NAME"
$ws.Range("C1").Value = "QUANTITY"
$ws.Range("D1").Value = "COST/UNIT"
$ws.Range("E1").Value = "TOTAL COST"
$ws.Range("F1").Value = "URL"

$headerRange = $ws.Range("A1:F1")
$headerRange.Style = "Heading 1"
$headerRange.HorizontalAlignment = -4108
$headerRange.RowHeight = 20.25
$headerRange.Borders.Item(9).LineStyle = 1
$headerRange.Borders.Item(9).Weight = 4

# --- Row 2 ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Utilite Standard"
$ws.Range("C2").Value = 16
$ws.Range("D2").Value = 159
$ws.Range("E2").Value = 2544
$ws.Range("F2").Value = "http://www.compulab.co.il/utilite-computer/web/utilite-availability"

# --- Row 3 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "HP 2920-24G Network Switch"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 872.95
$ws.Range("E3").Value = 872.95
$ws.Range("F3").Value = "http://www.amazon.com/HP-J9726A-2920-24G-Switch/dp/B00BJ42JYG/ref=sr_1_1?ie=UTF8&qid=1422292309&sr=8-1&keywords=J9726A"

# --- Number formats ---
$ws.Range("D2").NumberFormat = '"$"#,##0_);[Red]("$"#,##0)'
$ws.Range("E2").NumberFormat = '"$"#,##0_);[Red]("$"#,##0)'
$ws.Range("D3:E3").NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'

# --- Center alignment for data rows ---
$ws.Range("A2:F3").HorizontalAlignment = -4108

# --- Hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("F2"), "http://www.compulab.co.il/utilite-computer/web/utilite-availability") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "http://www.amazon.com/HP-J9726A-2920-24G-Switch/dp/B00BJ42JYG/ref=sr_1_1?ie=UTF8&qid=1422292309&sr=8-1&keywords=J9726A") | Out-Null

# --- Theme color tweak (dk1 text colour) ---
$wb.Theme.ThemeColorScheme.Colors(1).RGB = 0x373737

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("D17").Select()
